$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "67.627.37"
$ws.Range("E2").Value = "  +0.81%  "
# Row 3
$ws.Range("D3").Value = "2.495.07"
$ws.Range("E3").Value = "  +0.86%  "
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.11%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.61"
$ws.Range("E5").Value = "  +0.64%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "175.94"
$ws.Range("E6").Value = "  +4.05%  "
# Row 7
$ws.Range("E7").Value = "  +0.03%  "
# Row 8
$ws.Range("E8").Value = "  +0.26%  "
# Row 9
$ws.Range("E9").Value = "  +3.90%  "
# Row 10
$ws.Range("E10").Value = "  +0.33%  "
# Row 11
$ws.Range("E11").Value = "  +2.22%  "
# Row 12
$ws.Range("E12").Value = "  +0.05%  "
# Row 13
$ws.Range("D13").Value = "2.937.11"
$ws.Range("E13").Value = "  +0.42%  "
# Row 14
$ws.Range("E14").Value = "  +0.61%  "
# Row 15
$ws.Range("D15").Value = "67.521.21"
$ws.Range("E15").Value = "  +0.65%  "
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000171"
$ws.Range("E16").Value = "  +1.44%  "
# Row 17
$ws.Range("D17").Value = "2.487.35"
$ws.Range("E17").Value = "  +0.67%  "
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.05"
$ws.Range("E18").Value = "  -0.27%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.47"
$ws.Range("E19").Value = "  -0.76%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "351.95"
$ws.Range("E20").Value = "  +0.02%  "
# Row 21
$ws.Range("E21").Value = "  -0.72%  "
# Row 22
$ws.Range("E22").Value = "  +0.14%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.50"
$ws.Range("E23").Value = "  +2.32%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.23"
$ws.Range("E24").Value = "  -0.49%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.80"
$ws.Range("E25").Value = "  -1.22%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.28"
$ws.Range("E26").Value = "  +0.82%  "
# Row 27
$ws.Range("D27").Value = "2.614.01"
$ws.Range("E27").Value = "  +0.69%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  +0.00%  "
# Row 29
$ws.Range("D29").Value = "0.0₃0912"
$ws.Range("E29").Value = "  +1.00%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "508.44"
$ws.Range("E30").Value = "  -0.28%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.84"
$ws.Range("E31").Value = "  +1.77%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.25"
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.78"
$ws.Range("E33").Value = "  +0.45%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("E34").Value = "  +0.02%  "
# Row 35
$ws.Range("E35").Value = "  +5.14%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "160.63"
$ws.Range("E36").Value = "  +0.81%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "18.70"
$ws.Range("E37").Value = "  +0.22%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.32"
$ws.Range("E38").Value = "  -0.05%  "
# Row 39
$ws.Range("E39").Value = "  +0.56%  "
# Row 40
$ws.Range("E40").Value = "  +0.06%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.70"
$ws.Range("E41").Value = "  +0.89%  "
# Row 42
$ws.Range("E42").Value = "  +0.63%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.87"
$ws.Range("E43").Value = "  +0.98%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.43"
$ws.Range("E44").Value = "  +2.20%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "143.65"
$ws.Range("E45").Value = "  +1.78%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.50"
$ws.Range("E46").Value = "  +1.66%  "
# Row 47
$ws.Range("D47").Value = "0.0₆0260"
$ws.Range("E47").Value = "  +2.04%  "
# Row 48
$ws.Range("E48").Value = "  -0.11%  "
# Row 49
$ws.Range("E49").Value = "  +1.80%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.58"
$ws.Range("E50").Value = "  -0.37%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.587"
$ws.Range("E51").Value = "  +0.96%  "
